$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 21666.334
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 21666.334
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 21666.334
$ws.Range("N3").Value = -21894.334
$ws.Range("H102").Value = 21666.334
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 21666.334
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 21666.334
$ws.Range("N102").Value = -28156.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 950
$ws.Range("I25").Value = 950
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 950
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -548
$ws.Range("H32").Value = 15300.538
$ws.Range("I32").Value = 13242.333
$ws.Range("J32").Value = 39999
$ws.Range("K32").Value = 13242.333
$ws.Range("L32").Value = 39999
$ws.Range("M32").Value = -12955.333
$ws.Range("N32").Value = -40573
$ws.Range("H96").Value = 35172
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 35172
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 35172
$ws.Range("N96").Value = -40664
$ws.Range("H102").Value = 1850
$ws.Range("I102").Value = 1958.3334
$ws.Range("J102").Value = 1200
$ws.Range("K102").Value = 1958.3334
$ws.Range("L102").Value = 1200
$ws.Range("M102").Value = -336.3334
$ws.Range("N102").Value = -4444
$ws.Range("H112").Value = 77388.5
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 77388.5
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 77388.5
$ws.Range("N112").Value = -80342.5
$ws.Range("H132").Value = 3693.625
$ws.Range("I132").Value = 3693.625
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11080.875
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8550.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 500
$ws.Range("I12").Value = 500
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 500
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -332
$ws.Range("H80").Value = 1324.6666
$ws.Range("I80").Value = 777.5
$ws.Range("J80").Value = 1598.25
$ws.Range("K80").Value = 777.5
$ws.Range("L80").Value = 1598.25
$ws.Range("M80").Value = 220.5
$ws.Range("N80").Value = -3594.25
$ws.Range("H83").Value = 1324.6666
$ws.Range("I83").Value = 777.5
$ws.Range("J83").Value = 1598.25
$ws.Range("K83").Value = 3887.5
$ws.Range("L83").Value = 7991.25
$ws.Range("M83").Value = 1104.5
$ws.Range("N83").Value = -17975.25
$ws.Range("H94").Value = 5005
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 5005
$ws.Range("K94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("M94").Value = 5005
$ws.Range("N94").Value = -5907
$ws.Range("H99").Value = 1122.5
$ws.Range("I99").Value = 995
$ws.Range("J99").Value = 1165
$ws.Range("K99").Value = 995
$ws.Range("L99").Value = 1165
$ws.Range("M99").Value = 503
$ws.Range("N99").Value = -4161
$ws.Range("H134").Value = 1239.8
$ws.Range("I134").Value = 1239.8
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3719.4
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1184.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3757.625
$ws.Range("I31").Value = 2639.6
$ws.Range("J31").Value = 5621
$ws.Range("K31").Value = 2639.6
$ws.Range("L31").Value = 5621
$ws.Range("M31").Value = -2344.6
$ws.Range("N31").Value = -6211
$ws.Range("H34").Value = 3757.625
$ws.Range("I34").Value = 2639.6
$ws.Range("J34").Value = 5621
$ws.Range("K34").Value = 2639.6
$ws.Range("L34").Value = 5621
$ws.Range("M34").Value = -2437.6
$ws.Range("N34").Value = -6025
$ws.Range("H58").Value = 4453
$ws.Range("I58").Value = 3270.6667
$ws.Range("J58").Value = 8000
$ws.Range("K58").Value = 3270.6667
$ws.Range("L58").Value = 8000
$ws.Range("M58").Value = -3067.6667
$ws.Range("N58").Value = -8406
$ws.Range("H103").Value = 29903.334
$ws.Range("I103").Value = 29903.334
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 29903.334
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -28731.334
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").ClearContents()
$ws.Range("N108").Value = 0
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H136").Value = 4453
$ws.Range("I136").Value = 3270.6667
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 9812.000100000001
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -7262.000100000001
$ws.Range("N136").Value = -29100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1108
$ws.Range("I4").Value = 1108
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3324
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -3212
$ws.Range("H38").Value = 4551
$ws.Range("I38").Value = 9000
$ws.Range("J38").Value = 102
$ws.Range("K38").Value = 27000
$ws.Range("L38").Value = 306
$ws.Range("M38").Value = -26653
$ws.Range("N38").Value = -1000
$ws.Range("H68").Value = 714.2857
$ws.Range("I68").Value = 500
$ws.Range("J68").Value = 1000
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -689
$ws.Range("N68").Value = -4622
$ws.Range("H71").Value = 714.2857
$ws.Range("I71").Value = 500
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 4500
$ws.Range("L71").Value = 9000
$ws.Range("M71").Value = -444
$ws.Range("N71").Value = -17112
$ws.Range("H80").Value = 908.75
$ws.Range("I80").Value = 896.6667
$ws.Range("J80").Value = 945
$ws.Range("K80").Value = 2690.0001
$ws.Range("L80").Value = 2835
$ws.Range("M80").Value = -1754.0001
$ws.Range("N80").Value = -4707
$ws.Range("H83").Value = 908.75
$ws.Range("I83").Value = 896.6667
$ws.Range("J83").Value = 945
$ws.Range("K83").Value = 8070.0003
$ws.Range("L83").Value = 8505
$ws.Range("M83").Value = -3390.0003
$ws.Range("N83").Value = -17865
$ws.Range("H131").Value = 2511.7273
$ws.Range("I131").Value = 1604.8334
$ws.Range("J131").Value = 3600
$ws.Range("K131").Value = 4814.5002
$ws.Range("L131").Value = 10800
$ws.Range("M131").Value = 225.4997999999996
$ws.Range("N131").Value = -20880

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 7500
$ws.Range("I4").Value = 7500
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 7500
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -7388
$ws.Range("H43").Value = 3512.75
$ws.Range("I43").Value = 3512.75
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 3512.75
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -3361.75
$ws.Range("H118").Value = 99310
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 99310
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 99310
$ws.Range("N118").Value = -102624

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1941.5
$ws.Range("I61").Value = 1888
$ws.Range("J61").Value = 1995
$ws.Range("K61").Value = 1888
$ws.Range("L61").Value = 1995
$ws.Range("M61").Value = -1686
$ws.Range("N61").Value = -2399
$ws.Range("H113").Value = 1941.5
$ws.Range("I113").Value = 1888
$ws.Range("J113").Value = 1995
$ws.Range("K113").Value = 1888
$ws.Range("L113").Value = 1995
$ws.Range("M113").Value = 282
$ws.Range("N113").Value = -6335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2310
$ws.Range("I132").Value = 1887.75
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 5663.25
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -3133.25
$ws.Range("N132").Value = -17057
